$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value2 = "Baggy Koyu Mavi"
$ws.Range("B7").Value2 = "300 Tl"
$ws.Range("C7").Value2 = "Jeans"
$ws.Range("D7").Value2 = "BAG5.jpg"
$ws.Range("E7").Value2 = $ws.Range("E6").Value2
$ws.Range("F7").Value2 = "Var"

$ws.Range("E12").Select()
